# Actualización automática 2025-08-29 14:00:09
#
# Updates the August ("agosto") sales figures for two asesor/cliente
# transactions plus one additional transaction, and refreshes every value
# that is derived from them across the three report sheets:
#   - "VENTAS POR GRUPO"     : per-client sales broken down by product group
#   - "VENTA MENSUAL"        : per-client sales broken down by month
#   - "CUMPLIMIENTO MENSUAL" : per-group totals vs. budget (PRESUPUESTO)
#
# The source workbook stores everything as plain cached values (no
# formulas), so the "recalculation" has to be written back explicitly.

$wb = $excel.ActiveWorkbook

$wsGrupo  = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl  = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- "VENTAS POR GRUPO": raw per-client, per-group figures ---------------
# Row 11 = BUELE MORENO JOSELITO ROYERS, column P = NO RESURTIBLES
$wsGrupo.Range("P11").Value = 1706.16
# Row 13 = CABRERA VALVERDE ANAHI FERNANDA, column M = PORCELANATO
$wsGrupo.Range("M13").Value = 3844.25
# Row 52 = WONG SANCHEZ CLAUDIA PAULINA, column O = SAL SOLUBLE
$wsGrupo.Range("O52").Value = 2917

# --- "VENTA MENSUAL": same transactions, column F = agosto ----------------
$wsMensual.Range("F11").Value = 1706.16
$wsMensual.Range("F13").Value = 3199.36
$wsMensual.Range("F52").Value = 4750.73
# Row 55 = TOTAL column (sum of agosto for all clients)
$wsMensual.Range("F55").Value = 110299.61

# --- "CUMPLIMIENTO MENSUAL": PRESUPUESTO / VENTA / POR CUMPLIR / CUMPLIMIENTO
# Row 10 = NO RESURTIBLES group total
$wsCumpl.Range("D10").Value = 1746.49
$wsCumpl.Range("E10").Value = -445.99
$wsCumpl.Range("F10").Value = 1.342937331795463

# Row 16 = PORCELANATO group total
$wsCumpl.Range("D16").Value = 56227.15
$wsCumpl.Range("E16").Value = -167.4500000000044
$wsCumpl.Range("F16").Value = 1.002986994222231

# Row 18 = SAL SOLUBLE group total
$wsCumpl.Range("D18").Value = 13661.64
$wsCumpl.Range("E18").Value = -10461.64
$wsCumpl.Range("F18").Value = 4.2692625

# Row 19 = TOTAL row
$wsCumpl.Range("D19").Value = 110299.61
$wsCumpl.Range("E19").Value = 7140.080645179147
$wsCumpl.Range("F19").Value = 0.9392021504318203
